$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '63.915.50'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  -3.13%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.152.27'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  -8.49%  '

$ws.Range("E4").Value = '  +0.06%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '555.81'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -4.84%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '169.85'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -1.97%  '

$ws.Range("E7").Value = '  +0.04%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.603'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  -0.08%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '3.154.89'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -8.43%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.122'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -6.49%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.56'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -5.08%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.393'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  -4.05%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '3.710.44'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  -8.28%  '

$ws.Range("E14").Value = '  +0.12%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '27.25'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  -5.64%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '64.060.30'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  -3.00%  '

$ws.Range("E17").Value = '  -5.65%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.163.79'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  -8.39%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '5.63'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -5.19%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.94'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -6.51%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '350.31'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -4.80%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '7.11'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -7.23%  '

$ws.Range("E24").Value = '  -5.43%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.498'
$ws.Range("D25").ClearFormats()

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.0000116'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -4.04%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.35'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -4.23%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.174'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -1.62%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.999'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  +0.07%  '

$ws.Range("E30").Value = '  -0.03%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '5.55'
$ws.Range("D31").ClearFormats()

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.87'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -5.41%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '21.90'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -7.27%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.56'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -6.59%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.18'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -9.14%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '156.68'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -2.52%  '

$ws.Range("E37").Value = '  -7.00%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.798'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -9.38%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '25.87'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -10.33%  '

$ws.Range("E40").Value = '  -5.02%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.66'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -6.08%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.628.34'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -4.73%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '4.12'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -7.51%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '5.97'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -7.72%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0646'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -5.05%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '38.56'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -3.96%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '323.32'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -0.55%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '23.26'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -4.35%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0268'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -7.45%  '

$ws.Range("E50").Value = '  -0.77%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.00'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +0.02%  '
